# Decode-password edit:
#  - Password cells E2:E5 on the DATA sheet now store the Base64-encoded
#    password ("VGVzdEAxMjM=" == base64("Test@123")) instead of the
#    plaintext value, so downstream code can add a decoder step.
#  - The five individual mailto hyperlinks that used to sit on E2,E3,E4,E5
#    (all pointing at the old plaintext password) are replaced by a single
#    hyperlink spanning the merged-looking E2:E5 selection; K4's hyperlink
#    is kept.
#  - The DATA sheet becomes the active/selected sheet with E2:E5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# --- Rebuild hyperlinks ---------------------------------------------------
# Remove every existing hyperlink on the sheet (the old per-cell mailto
# links on E2:E5) and recreate the one that should remain (K4) plus the
# new merged link covering E2:E5. Do this BEFORE touching cell values,
# since Hyperlinks.Add() stamps its TextToDisplay into the anchor range's
# first cell.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K4"), "mailto:test@ram.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2:E5"), "mailto:test@ram.com", "", "", "test@ram.com") | Out-Null

# --- Update the password cells: plaintext -> Base64-encoded value -------
$encodedPassword = "'VGVzdEAxMjM="

$eCells = @("E2", "E3", "E4", "E5")
foreach ($addr in $eCells) {
    $cell = $ws.Range($addr)
    $cell.Style = "Hyperlink"
    $cell.Value = $encodedPassword
}

# --- Sheet/selection activation ------------------------------------------
# The DATA sheet becomes the active tab, with E2:E5 as the selection
# (keeping the existing topLeftCell scroll position).
$ws.Activate()
$ws.Range("E2:E5").Select()
